$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.799.53'
$ws.Range("E2").Value = '  +0.02%  '

# Row 3
$ws.Range("D3").Value = '3.537.75'
$ws.Range("E3").Value = '  +1.27%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.61'
$ws.Range("E5").Value = '  -0.47%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '194.80'
$ws.Range("E6").Value = '  +0.07%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  -0.34%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.204'
$ws.Range("E9").Value = '  -4.12%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.647'
$ws.Range("E10").Value = '  -1.71%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.39'
$ws.Range("E11").Value = '  -0.22%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000303'
$ws.Range("E12").Value = '  -1.00%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.50'
$ws.Range("E13").Value = '  -1.14%  '

# Row 14
$ws.Range("D14").Value = '4.092.42'

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '593.07'
$ws.Range("E15").Value = '  -1.36%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '12.78'
$ws.Range("E16").Value = '  +1.06%  '

# Row 17
$ws.Range("D17").Value = '69.896.15'
$ws.Range("E17").Value = '  +0.05%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.03'
$ws.Range("E18").Value = '  +0.66%  '

# Row 19
$ws.Range("D19").Value = '3.538.02'
$ws.Range("E19").Value = '  +1.15%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.123'
$ws.Range("E20").Value = '  +1.81%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.983'
$ws.Range("E21").Value = '  -0.55%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.82'
$ws.Range("E22").Value = '  -0.37%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '103.11'
$ws.Range("E23").Value = '  -1.87%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.12'
$ws.Range("E24").Value = '  +1.12%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.63'
$ws.Range("E25").Value = '  -0.23%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.04'
$ws.Range("E26").Value = '  -1.16%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.73'
$ws.Range("E27").Value = '  -2.03%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.52'
$ws.Range("E28").Value = '  -3.12%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.17'
$ws.Range("E29").Value = '  -2.41%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.03'
$ws.Range("E30").Value = '  -1.98%  '

# Row 31
$ws.Range("E31").Value = '  -4.61%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.31'
$ws.Range("E32").Value = '  -2.72%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.115'
$ws.Range("E33").Value = '  -0.20%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.42'
$ws.Range("E34").Value = '  -1.09%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.22'
$ws.Range("E35").Value = '  +5.93%  '

# Row 36
$ws.Range("D36").Value = '3.826.64'
$ws.Range("E36").Value = '  +3.84%  '

# Row 37
$ws.Range("D37").Value = '0.0₃0822'
$ws.Range("E37").Value = '  +4.37%  '

# Row 38
$ws.Range("E38").Value = '  +0.24%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '511.36'
$ws.Range("E39").Value = '  -1.59%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.391'
$ws.Range("E40").Value = '  +0.12%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.56'
$ws.Range("E41").Value = '  -0.25%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '36.51'
$ws.Range("E42").Value = '  -0.62%  '

# Row 43
$ws.Range("E43").Value = '  -2.42%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0447'
$ws.Range("E44").Value = '  -2.94%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.31'
$ws.Range("E45").Value = '  -0.15%  '

# Row 46
$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.81'
$ws.Range("E46").Value = '  -1.20%  '

# Row 47
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.139'
$ws.Range("E47").Value = '  -1.00%  '

# Row 48
$ws.Range("E48").Value = '  -0.02%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.50'
$ws.Range("E49").Value = '  -3.14%  '

# Row 50
$ws.Range("E50").Value = '  +3.17%  '

# Row 51
$ws.Range("E51").Value = '  +2.17%  '
